# Append new Shipping Confirmation rows (order-id, order-item-id, carrier-code, tracking-number)
# to the "ShippingConfirmation" sheet, rows 18-44, matching columns A, B, E, G
# (order-id, order-item-id, carrier-code, tracking-number) used by the existing rows 2-17.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ShippingConfirmation")

# Each entry: order-id, order-item-id, carrier-code, tracking-number
$newRows = @(
  @("112-3616143-6610663","14-12923-67816","UPS","1Z8126VRYW33280602"),
  @("112-2803937-5654632","702914437","UPS","1Z8126VRYW33352374"),
  @("113-6413731-1090631","702942056","USPS","9234690185108900016648"),
  @("111-5534447-2556238","02-12942-39922","UPS","1ZA5T352YW01399846"),
  @("112-1404559-5127409","07-12935-80266","UPS","1ZA5T352YW01398936"),
  @("111-9945324-9853836","703212138","UPS","1ZA5T3520201400315"),
  @("114-8327013-2852237","10-12931-62213","UPS","1ZA5T352YW01397384"),
  @("114-2733954-5928252","703195679","UPS","1ZA5T3520201399944"),
  @("111-1960000-7331436","703042547","UPS","1ZA5T352YW01397982"),
  @("113-1778148-4749838","01-12947-15133","UPS","1Z3024W2YW00802738"),
  @("114-3140840-2805024","15-12928-93200","UPS","1ZA5T3520301443287"),
  @("112-0163782-4789872","18-12925-65132","UPS","1Z3024W20300799170"),
  @("114-4592665-1282659","14-12930-44483","UPS","1ZA5T352YW01446875"),
  @("113-6596600-2407435","10-12935-74847","UPS","1Z8126VRYW33419534"),
  @("111-0287793-5740220","19-12924-17640","UPS","1ZA5T352YW01447463"),
  @("113-4054909-3593014","26-12915-31709","UPS","1ZA5T3520301443241"),
  @("113-9429656-1739405","08-12938-69957","UPS","1Z3024W2YW00806645"),
  @("112-7517354-7596220","19-12924-40995","UPS","1ZA5T352YW01448944"),
  @("113-1204470-2667418","19-12924-16048","UPS","1ZA5T352YW01446740"),
  @("114-3974135-4995407","27-12914-97960","UPS","1Z3024W2YW00804941"),
  @("114-9444413-5227407","704067812","UPS","1Z3024W20300803682"),
  @("112-6419306-5168246","25-12917-96224","UPS","1Z3024W2YW00806243"),
  @("113-9205704-4961051","10-12936-21810","UPS","1Z3024W2YW00805842"),
  @("113-2299668-6301016","25-12918-29177","UPS","1Z8126VRYW33426768"),
  @("111-4714273-6097824","703999666","UPS","1Z8126VRYW33415234"),
  @("112-0513920-0415415","08-12938-29265","UPS","1Z8126VRYW33421656"),
  @("112-3978083-8357811","704073723","UPS","1Z3024W20300806769")
)

$row = 18
foreach ($item in $newRows) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 1).Style = "Normal"

    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 2).Style = "Normal"

    $ws.Cells.Item($row, 5).Value = $item[2]
    $ws.Cells.Item($row, 5).Style = "Normal"

    $ws.Cells.Item($row, 7).Value = $item[3]
    $ws.Cells.Item($row, 7).Style = "Normal"

    $row = $row + 1
}

# Reflect the new extent of used data in the view (selection over rows 2:44, same as
# what a user would see after selecting/filling down through the newly added rows).
$ws.Rows("2:44").Select()
